# ---------------------------------------------------------------------------
# prog1_data.xlsx — "added more values to spreadsheet"
#
# For each of the n=2 / n=3 / n=4 experiment blocks, the sheet originally
# only had a header row plus an empty "Average MST Weight" row (no data had
# been recorded for them yet). This edit fills in the "Average MST Weight"
# row and adds three more rows per block: "Max included edge", "(trials)"
# and the k(n) formula-description label -- mirroring the fully populated
# n=0 block at the top of the sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows after each "Average MST Weight" row, working bottom-up so
# earlier anchors keep their original row numbers while we work.
$ws.Rows("18:20").Insert()
$ws.Rows("15:17").Insert()
$ws.Rows("12:14").Insert()

# --- n=2 block: fill Average MST Weight (row 11) + new rows 12-14 ---
$ws.Range("E11").Value = 2.723387
$ws.Range("F11").Value = 3.869813
$ws.Range("G11").Value = 5.437021
$ws.Range("H11").Value = 7.611348
$ws.Range("I11").Value = 10.660027
$ws.Range("J11").Value = 14.981488
$ws.Range("K11").Value = 21.055611
$ws.Range("D12").Value = "Max included edge"
$ws.Range("E12").Value = 0.695356
$ws.Range("F12").Value = 0.448285
$ws.Range("G12").Value = 0.352937
$ws.Range("H12").Value = 0.240012
$ws.Range("I12").Value = 0.184145
$ws.Range("J12").Value = 0.130891
$ws.Range("K12").Value = 0.082371
$ws.Range("D13").Value = "(trials)"
$ws.Range("E13").Value = 1000
$ws.Range("F13").Value = 1000
$ws.Range("G13").Value = 1000
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 100
$ws.Range("D14").Value = "0.025+(1.25)*4.53261/(n^0.754872)"

# --- n=3 block: fill Average MST Weight (row 17) + new rows 18-20 ---
$ws.Range("E17").Value = 4.523913
$ws.Range("F17").Value = 7.162139
$ws.Range("G17").Value = 11.240602
$ws.Range("H17").Value = 17.631735
$ws.Range("I17").Value = 27.600925
$ws.Range("J17").Value = 43.319759
$ws.Range("K17").Value = 68.157043
$ws.Range("D18").Value = "Max included edge"
$ws.Range("E18").Value = 0.82371
$ws.Range("F18").Value = 0.727282
$ws.Range("G18").Value = 0.533838
$ws.Range("H18").Value = 0.43611
$ws.Range("I18").Value = 0.356399
$ws.Range("J18").Value = 0.25657
$ws.Range("K18").Value = 0.210055
$ws.Range("D19").Value = "(trials)"
$ws.Range("E19").Value = 1000
$ws.Range("F19").Value = 1000
$ws.Range("G19").Value = 1000
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 100
$ws.Range("D20").Value = "0.025+(1.25)*4.53261/(n^0.754872)"

# --- n=4 block: fill Average MST Weight (row 23) + new rows 24-26 ---
$ws.Range("E23").Value = 6.13595
$ws.Range("F23").Value = 10.339508
$ws.Range("G23").Value = 17.126459
$ws.Range("H23").Value = 28.387136
$ws.Range("I23").Value = 47.187328
$ws.Range("J23").Value = 78.174957
$ws.Range("K23").Value = 129.935715
$ws.Range("D24").Value = "Max included edge"
$ws.Range("E24").Value = 0.886194
$ws.Range("F24").Value = 0.736967
$ws.Range("G24").Value = 0.680934
$ws.Range("H24").Value = 0.59279
$ws.Range("I24").Value = 0.465153
$ws.Range("J24").Value = 0.384114
$ws.Range("K24").Value = 0.311332
$ws.Range("D25").Value = "(trials)"
$ws.Range("E25").Value = 1000
$ws.Range("F25").Value = 1000
$ws.Range("G25").Value = 1000
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 100
$ws.Range("D26").Value = "0.025+(1.25)*4.53261/(n^0.754872)"


# Column width tweaks (C got narrower, D got wider to fit the new labels)
$ws.Columns("C").ColumnWidth = 7.666666666666667
$ws.Columns("D").ColumnWidth = 30.830729166666668

# Leave the selection where the user ended up after entering the new data
$ws.Range("H29").Select()
